$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("A2").Value = "Day 1 - 23 h 42"
$ws.Range("C2").Value = "Day 2 - 3 h 54"
$ws.Range("A3").Value = "Day 2 - 3 h 54"
$ws.Range("C3").Value = "Day 2 - 8 h 6"
$ws.Range("A4").Value = "Day 2 - 11 h 54"
$ws.Range("C4").Value = "Day 2 - 15 h 24"
$ws.Range("A5").Value = "Day 2 - 20 h 0"
$ws.Range("C5").Value = "Day 2 - 23 h 30"
$ws.Range("A6").Value = "Day 3 - 0 h 24"
$ws.Range("C6").Value = "Day 3 - 3 h 54"
$ws.Range("A7").Value = "Day 3 - 3 h 54"
$ws.Range("C7").Value = "Day 3 - 7 h 24"
$ws.Range("A8").Value = "Day 4 - 3 h 54"
$ws.Range("C8").Value = "Day 4 - 8 h 6"
$ws.Range("A9").Value = "Day 5 - 4 h 0"
$ws.Range("C9").Value = "Day 5 - 8 h 12"
$ws = $wb.Worksheets.Item(2)
$ws.Range("A2").Value = "Day 2 - 2 h 54"
$ws.Range("C2").Value = "Day 2 - 7 h 54"
$ws.Range("A3").Value = "Day 2 - 7 h 54"
$ws.Range("C3").Value = "Day 2 - 12 h 54"
$ws.Range("A4").Value = "Day 2 - 23 h 54"
$ws.Range("C4").Value = "Day 3 - 4 h 54"
$ws.Range("A5").Value = "Day 4 - 0 h 0"
$ws.Range("C5").Value = "Day 4 - 5 h 0"
$ws.Range("A6").Value = "Day 4 - 7 h 54"
$ws.Range("C6").Value = "Day 4 - 11 h 18"
$ws.Range("A7").Value = "Day 4 - 20 h 0"
$ws.Range("C7").Value = "Day 4 - 23 h 24"
$ws.Range("A8").Value = "Day 4 - 23 h 54"
$ws.Range("C8").Value = "Day 5 - 4 h 18"
$ws.Range("A9").Value = "Day 5 - 4 h 18"
$ws.Range("C9").Value = "Day 5 - 8 h 42"
$ws = $wb.Worksheets.Item(3)
$ws.Range("A2").Value = "Day 2 - 3 h 54"
$ws.Range("C2").Value = "Day 2 - 7 h 18"
$ws.Range("A3").Value = "Day 2 - 20 h 0"
$ws.Range("C3").Value = "Day 2 - 23 h 24"
$ws.Range("A4").Value = "Day 2 - 23 h 30"
$ws.Range("C4").Value = "Day 3 - 3 h 54"
$ws.Range("A5").Value = "Day 3 - 3 h 54"
$ws.Range("C5").Value = "Day 3 - 8 h 18"
$ws.Range("A6").Value = "Day 3 - 21 h 42"
$ws.Range("C6").Value = "Day 4 - 1 h 6"
$ws.Range("A7").Value = "Day 4 - 1 h 6"
$ws.Range("C7").Value = "Day 4 - 4 h 30"
$ws.Range("A8").Value = "Day 4 - 4 h 30"
$ws.Range("C8").Value = "Day 4 - 7 h 54"
$ws.Range("A9").Value = "Day 4 - 7 h 54"
$ws.Range("C9").Value = "Day 4 - 11 h 18"
$ws.Range("A10").Value = "Day 4 - 20 h 48"
$ws.Range("C10").Value = "Day 4 - 23 h 54"
$ws.Range("A11").Value = "Day 4 - 23 h 54"
$ws.Range("C11").Value = "Day 5 - 3 h 0"
$ws.Range("A12").Value = "Day 5 - 4 h 12"
$ws.Range("C12").Value = "Day 5 - 7 h 54"
$ws.Range("A13").Value = "Day 5 - 7 h 54"
$ws.Range("C13").Value = "Day 5 - 11 h 36"
$ws = $wb.Worksheets.Item(4)
$ws.Range("A2").Value = "Day 2 - 23 h 54"
$ws.Range("C2").Value = "Day 3 - 4 h 0"
$ws.Range("A3").Value = "Day 3 - 20 h 0"
$ws.Range("C3").Value = "Day 4 - 0 h 6"
$ws.Range("A4").Value = "Day 4 - 3 h 54"
$ws.Range("C4").Value = "Day 4 - 8 h 12"
$ws.Range("A5").Value = "Day 5 - 4 h 0"
$ws.Range("C5").Value = "Day 5 - 8 h 18"
$ws = $wb.Worksheets.Item(5)
$ws.Range("A2").Value = "Day 3 - 3 h 54"
$ws.Range("C2").Value = "Day 3 - 7 h 30"
$ws.Range("A3").Value = "Day 4 - 4 h 0"
$ws.Range("C3").Value = "Day 4 - 7 h 36"
